$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 1222
$ws.Range("F5").Value = 1399
$ws.Range("G5").Value = 88
$ws.Range("F6").Value = 1723
$ws.Range("F7").Value = 6250
$ws.Range("F9").Value = 1843
$ws.Range("F10").Value = 489
$ws.Range("F16").Value = 7027
$ws.Range("F19").Value = 173
$ws.Range("F20").Value = 105
$ws.Range("F21").Value = 1717
$ws.Range("F26").Value = 1626
$ws.Range("F27").Value = 773
$ws.Range("F28").Value = 328
$ws.Range("F29").Value = 2
$ws.Range("F31").Value = 61
$ws.Range("F33").Value = 3899

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 353
$ws.Range("F8").Value = 442

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F5").Value = 250

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F5").Value = 1222
$ws.Range("F9").Value = 353
$ws.Range("F10").Value = 1399
$ws.Range("G10").Value = 88
$ws.Range("F11").Value = 250
$ws.Range("F12").Value = 1723
$ws.Range("F13").Value = 6250
$ws.Range("F14").Value = 1843
$ws.Range("F17").Value = 489
$ws.Range("F24").Value = 7027
$ws.Range("F27").Value = 174
$ws.Range("F28").Value = 105
$ws.Range("F29").Value = 1717
$ws.Range("F33").Value = 1626
$ws.Range("F34").Value = 773
$ws.Range("F36").Value = 328
$ws.Range("F37").Value = 2
$ws.Range("F44").Value = 3899
